$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: content sourced from original row 46
$ws.Range("A39").Value2 = 111880601
$ws.Range("B39").Value2 = 88966
$ws.Range("E39").Value2 = 5754
$ws.Range("F39").Value2 = "Gultoppig fingersvamp"
$ws.Range("G39").Value2 = "Ramaria testaceoflava"
$ws.Range("H39").Value2 = "(Bres.) Corner"
$ws.Range("I39").Value = "'4"
$ws.Range("Q39").Value2 = 509942
$ws.Range("R39").Value2 = 6753225
$ws.Range("AJ39").Value2 = "tall"
$ws.Range("AK39").Value2 = "Pinus sylvestris"
$ws.Range("AL39").ClearContents()
$ws.Range("AO39").Value2 = "Pinus sylvestris"
$ws.Range("Z39").ClearContents()
$ws.Range("AB39").ClearContents()

# Row 40: content sourced from original row 42
$ws.Range("A40").Value2 = 111880475
$ws.Range("B40").Value2 = 88966
$ws.Range("E40").Value2 = 5754
$ws.Range("F40").Value2 = "Gultoppig fingersvamp"
$ws.Range("G40").Value2 = "Ramaria testaceoflava"
$ws.Range("H40").Value2 = "(Bres.) Corner"
$ws.Range("I40").Value = "'2"
$ws.Range("Q40").Value2 = 509958
$ws.Range("R40").Value2 = 6753363
$ws.Range("AJ40").Value2 = "gran"
$ws.Range("AK40").Value2 = "Picea abies"
$ws.Range("AL40").ClearContents()
$ws.Range("AO40").Value2 = "Picea abies"
$ws.Range("Z40").ClearContents()
$ws.Range("AB40").ClearContents()

# Row 41: content sourced from original row 44
$ws.Range("A41").Value2 = 111880500
$ws.Range("B41").Value2 = 88966
$ws.Range("E41").Value2 = 5754
$ws.Range("F41").Value2 = "Gultoppig fingersvamp"
$ws.Range("G41").Value2 = "Ramaria testaceoflava"
$ws.Range("H41").Value2 = "(Bres.) Corner"
$ws.Range("I41").Value = "'4"
$ws.Range("Q41").Value2 = 509899
$ws.Range("R41").Value2 = 6753571
$ws.Range("AJ41").Value2 = "gran"
$ws.Range("AK41").Value2 = "Picea abies"
$ws.Range("AL41").ClearContents()
$ws.Range("AO41").Value2 = "Picea abies"
$ws.Range("Z41").ClearContents()
$ws.Range("AB41").ClearContents()

# Row 42: content sourced from original row 45
$ws.Range("A42").Value2 = 111880484
$ws.Range("B42").Value2 = 90658
$ws.Range("E42").Value2 = 4361
$ws.Range("F42").Value2 = "Orange taggsvamp"
$ws.Range("G42").Value2 = "Hydnellum aurantiacum"
$ws.Range("H42").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I42").Value = "'11"
$ws.Range("Q42").Value2 = 509901
$ws.Range("R42").Value2 = 6753525
$ws.Range("AJ42").Value2 = "tall"
$ws.Range("AK42").Value2 = "Pinus sylvestris"
$ws.Range("AL42").ClearContents()
$ws.Range("AO42").Value2 = "Pinus sylvestris"
$ws.Range("Z42").ClearContents()
$ws.Range("AB42").ClearContents()

# Row 43: content sourced from original row 48
$ws.Range("A43").Value2 = 111880574
$ws.Range("B43").Value2 = 90658
$ws.Range("E43").Value2 = 4361
$ws.Range("F43").Value2 = "Orange taggsvamp"
$ws.Range("G43").Value2 = "Hydnellum aurantiacum"
$ws.Range("H43").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I43").Value = "'2"
$ws.Range("Q43").Value2 = 509596
$ws.Range("R43").Value2 = 6753392
$ws.Range("AJ43").Value2 = "tall"
$ws.Range("AK43").Value2 = "Pinus sylvestris"
$ws.Range("AL43").ClearContents()
$ws.Range("AO43").Value2 = "Pinus sylvestris"
$ws.Range("Z43").ClearContents()
$ws.Range("AB43").ClearContents()

# Row 44: content sourced from original row 41
$ws.Range("A44").Value2 = 111880532
$ws.Range("B44").Value2 = 90652
$ws.Range("E44").Value2 = 3100
$ws.Range("F44").Value2 = "Talltaggsvamp"
$ws.Range("G44").Value2 = "Bankera fuligineoalba"
$ws.Range("H44").Value2 = "(Schmidt : Fr.) Pouzar"
$ws.Range("I44").Value = "'2"
$ws.Range("Q44").Value2 = 509683
$ws.Range("R44").Value2 = 6753541
$ws.Range("AJ44").Value2 = "tall"
$ws.Range("AK44").Value2 = "Pinus sylvestris"
$ws.Range("AL44").ClearContents()
$ws.Range("AO44").Value2 = "Pinus sylvestris"
$ws.Range("Z44").ClearContents()
$ws.Range("AB44").ClearContents()

# Row 45: content sourced from original row 40
$ws.Range("A45").Value2 = 111880580
$ws.Range("B45").Value2 = 90658
$ws.Range("E45").Value2 = 4361
$ws.Range("F45").Value2 = "Orange taggsvamp"
$ws.Range("G45").Value2 = "Hydnellum aurantiacum"
$ws.Range("H45").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I45").Value = "'3"
$ws.Range("Q45").Value2 = 509755
$ws.Range("R45").Value2 = 6753236
$ws.Range("AJ45").Value2 = "tall"
$ws.Range("AK45").Value2 = "Pinus sylvestris"
$ws.Range("AL45").ClearContents()
$ws.Range("AO45").Value2 = "Pinus sylvestris"
$ws.Range("Z45").ClearContents()
$ws.Range("AB45").ClearContents()

# Row 46: content sourced from original row 39
$ws.Range("A46").Value2 = 111880462
$ws.Range("B46").Value2 = 88966
$ws.Range("E46").Value2 = 5754
$ws.Range("F46").Value2 = "Gultoppig fingersvamp"
$ws.Range("G46").Value2 = "Ramaria testaceoflava"
$ws.Range("H46").Value2 = "(Bres.) Corner"
$ws.Range("I46").Value = "'1"
$ws.Range("Q46").Value2 = 509970
$ws.Range("R46").Value2 = 6753250
$ws.Range("AJ46").Value2 = "tall"
$ws.Range("AK46").Value2 = "Pinus sylvestris"
$ws.Range("AL46").Value2 = "vid tallar"
$ws.Range("AO46").Value2 = "Pinus sylvestris # vid tallar"
$ws.Range("Z46").ClearContents()
$ws.Range("AB46").ClearContents()

# Row 47: content sourced from original row 47
$ws.Range("A47").Value2 = 111880562
$ws.Range("B47").Value2 = 90658
$ws.Range("E47").Value2 = 4361
$ws.Range("F47").Value2 = "Orange taggsvamp"
$ws.Range("G47").Value2 = "Hydnellum aurantiacum"
$ws.Range("H47").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I47").Value = "'3"
$ws.Range("Q47").Value2 = 509658
$ws.Range("R47").Value2 = 6753521
$ws.Range("AJ47").Value2 = "tall"
$ws.Range("AK47").Value2 = "Pinus sylvestris"
$ws.Range("AL47").ClearContents()
$ws.Range("AO47").Value2 = "Pinus sylvestris"
$ws.Range("Z47").ClearContents()
$ws.Range("AB47").ClearContents()

# Row 48: content sourced from original row 49
$ws.Range("A48").Value2 = 111880509
$ws.Range("B48").Value2 = 90652
$ws.Range("E48").Value2 = 3100
$ws.Range("F48").Value2 = "Talltaggsvamp"
$ws.Range("G48").Value2 = "Bankera fuligineoalba"
$ws.Range("H48").Value2 = "(Schmidt : Fr.) Pouzar"
$ws.Range("I48").Value = "'6"
$ws.Range("Q48").Value2 = 509834
$ws.Range("R48").Value2 = 6753644
$ws.Range("AJ48").Value2 = "tall"
$ws.Range("AK48").Value2 = "Pinus sylvestris"
$ws.Range("AL48").ClearContents()
$ws.Range("AO48").Value2 = "Pinus sylvestris"
$ws.Range("Z48").ClearContents()
$ws.Range("AB48").ClearContents()

# Row 49: content sourced from original row 43
$ws.Range("A49").Value2 = 111880591
$ws.Range("B49").Value2 = 90658
$ws.Range("E49").Value2 = 4361
$ws.Range("F49").Value2 = "Orange taggsvamp"
$ws.Range("G49").Value2 = "Hydnellum aurantiacum"
$ws.Range("H49").Value2 = "(Batsch:Fr.) P.Karst."
$ws.Range("I49").Value = "'8"
$ws.Range("Q49").Value2 = 509822
$ws.Range("R49").Value2 = 6753234
$ws.Range("AJ49").Value2 = "tall"
$ws.Range("AK49").Value2 = "Pinus sylvestris"
$ws.Range("AL49").ClearContents()
$ws.Range("AO49").Value2 = "Pinus sylvestris"
$ws.Range("Z49").ClearContents()
$ws.Range("AB49").ClearContents()
